# Generate Report for Handoff
# - Overview sheet: Status columns (zh-cn / de-de) move from "In Translation" to
#   "Ready for handoff", and the Latest Handoff Date bumps from 22:28:06 to 22:28:51.
# - zh-cn sheet: Latest Handoff Datetime bumps from 22:27:47 to 22:28:47.
# - de-de sheet: Latest Handoff Datetime bumps from 22:28:06 to 22:28:51.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet (row 2): B = zh-cn status, C = de-de status, D = Latest Handoff Date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-28-17 22:28:51"

# zh-cn sheet (row 2): C = Status, E = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-17 22:28:47"

# de-de sheet (row 2): C = Status, E = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-17 22:28:51"
